$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 503 ("「気分が良いぞう」I فيل GOOD ..." post), shifting all
# subsequent rows up by one.
$ws.Rows.Item(503).Delete()
